# Updates the Ky / Kz ("K values") columns (F and G) in the SSC_events
# table, per the commit "updated K values in spreadsheet".
#
# Table layout (row 1 = headers):
#   A=Date B=Time C=dBy D=dBz E=SYM-H Jump F=Ky G=Kz H=Usw
#   I=By,init J=By,final K=Bz,init L=Bz,final

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (event 2000-04-06)
$ws.Range("F2").Value = "2.4 x 10^-2 A/m"
$ws.Range("G2").Value = "1.0 x 10^-2 A/m"

# Row 3 (event 2001-11-24)
$ws.Range("F3").Value = "1.8 x 10^-2 A/m"
$ws.Range("G3").Value = "2.0 x 10^-2 A/m"

# Row 4 (event 2003-10-28)
$ws.Range("F4").Value = "-4.0 x 10^-3 A/m"
$ws.Range("G4").Value = "4.0 x 10^-3 A/m"

# Row 5 (event 2003-11-15, "ideal event" row 6 unaffected)
$ws.Range("F5").Value = "8.0 x 10^-3 A/m"
$ws.Range("G5").Value = "-4.8 x 10^-3 A/m"

# Cosmetic view state (best-effort; scroll position of the saved
# worksheet window moved so column E is the left-most visible column,
# with F5 remaining the active cell/selection).
$ws.Range("F5").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
